$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 201
$ws1.Range("F4").Value = 216
$ws1.Range("F6").Value = 1329
$ws1.Range("G6").Value = 61.2
$ws1.Range("F9").Value = 394
$ws1.Range("F10").Value = 449
$ws1.Range("F11").Value = 817
$ws1.Range("F12").Value = 221
$ws1.Range("F13").Value = 750
$ws1.Range("F14").Value = 318
$ws1.Range("F15").Value = 477
$ws1.Range("F17").Value = 1056
$ws1.Range("F18").Value = 501
$ws1.Range("F20").Value = 415
$ws1.Range("F21").Value = 107
$ws1.Range("F22").Value = 228
$ws1.Range("F26").Value = 451
$ws1.Range("F27").Value = 309

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 383
$ws2.Range("F5").Value = 52
$ws2.Range("F7").Value = 297
$ws2.Range("F11").Value = 161
$ws2.Range("F12").Value = 150

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 201
$ws4.Range("F6").Value = 216
$ws4.Range("F8").Value = 1329
$ws4.Range("G8").Value = 61.2
$ws4.Range("F11").Value = 383
$ws4.Range("F13").Value = 52
$ws4.Range("F14").Value = 394
$ws4.Range("F16").Value = 297
$ws4.Range("F17").Value = 449
$ws4.Range("F18").Value = 817
$ws4.Range("F19").Value = 221
$ws4.Range("F20").Value = 750
$ws4.Range("F21").Value = 318
$ws4.Range("F22").Value = 477
$ws4.Range("F24").Value = 1056
$ws4.Range("F25").Value = 501
$ws4.Range("F29").Value = 415
$ws4.Range("F31").Value = 107
$ws4.Range("F32").Value = 161
$ws4.Range("F33").Value = 228
$ws4.Range("F36").Value = 150
$ws4.Range("F41").Value = 451
$ws4.Range("F42").Value = 309
